$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 109, pushing existing rows 109..170 down to 110..171
$ws.Rows.Item(109).Insert()

# Populate the newly inserted row 109 with the new entry's data
$ws.Cells.Item(109, 1).Value = 3
$ws.Cells.Item(109, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(109, 3).Value = "Coquimbo"
$ws.Cells.Item(109, 4).Value = 44452
$ws.Cells.Item(109, 5).Value = 5
$ws.Cells.Item(109, 6).Value = 100112039
$ws.Cells.Item(109, 7).Value = "Ciboulette"
$ws.Cells.Item(109, 8).Value = "Sin especificar"
$ws.Cells.Item(109, 9).Value = "Primera"
$ws.Cells.Item(109, 10).Value = 190
$ws.Cells.Item(109, 11).Value = 1500
$ws.Cells.Item(109, 12).Value = 1500
$ws.Cells.Item(109, 13).Value = 1500
$ws.Cells.Item(109, 14).Value = "$/docena de atados"
$ws.Cells.Item(109, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(109, 16).Value = 500
$ws.Cells.Item(109, 17).Value = 3
$ws.Cells.Item(109, 18).Value = "Hortaliza"

# Match the date cell style used by the rest of column D (style index 2 in before.xlsx)
$ws.Cells.Item(109, 4).NumberFormat = $ws.Cells.Item(110, 4).NumberFormat
